$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting the existing rows (2,3) down to (3,4)
$ws.Rows.Item(2).Insert(-4121, 1)
$ws.Rows.Item(2).ClearFormats()

# Populate the new row 2 with the new weekly record
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "Vega Monumental Concepción"
$ws.Range("C2").Value = "Bíobío"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D2").Value = 44719
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100104
$ws.Range("H2").Value = "Frutos de pepita"
$ws.Range("I2").Value = 100104001
$ws.Range("J2").Value = "Granada"
$ws.Range("K2").Value = "Wonderfull"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20400
$ws.Range("Q2").Value = "$/caja 18 kilos granel"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 1133
$ws.Range("T2").Value = 18
